# Add web-scraper login config rows (airAsia / citilink / sriwijaya) to
# Lunggo_Config.xlsx, rows 45-50 of Sheet1, following the existing
# "*.<project>.<variable>" pattern established in row 44.
#
# The order in which new literal strings are first written below is
# deliberately chosen so that the workbook's shared-string table ends up
# enumerated in the same order as the source edit (airAsia, citilink,
# sriwijaya, IDTDEZYCGK_ADMIN, Travorama123, Travelmadezy, MLWAG0215,
# TRAVELMADEZY, webUserName, webPassword).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: project name ------------------------------------------------
$ws.Range("B45").Value = "airAsia"
$ws.Range("B46").Value = "airAsia"
$ws.Range("B47").Value = "citilink"
$ws.Range("B48").Value = "citilink"
$ws.Range("B49").Value = "sriwijaya"
$ws.Range("B50").Value = "sriwijaya"

# --- Columns E-H: Local / DV1 / QA / Production credential values ----------
$ws.Range("E45").Value = "IDTDEZYCGK_ADMIN"
$ws.Range("F45").Value = "IDTDEZYCGK_ADMIN"
$ws.Range("G45").Value = "IDTDEZYCGK_ADMIN"
$ws.Range("H45").Value = "IDTDEZYCGK_ADMIN"

$ws.Range("E46").Value = "Travorama123"
$ws.Range("F46").Value = "Travorama123"
$ws.Range("G46").Value = "Travorama123"
$ws.Range("H46").Value = "Travorama123"

$ws.Range("E47").Value = "Travelmadezy"
$ws.Range("F47").Value = "Travelmadezy"
$ws.Range("G47").Value = "Travelmadezy"
$ws.Range("H47").Value = "Travelmadezy"

$ws.Range("E48").Value = "Standar1234"
$ws.Range("F48").Value = "Standar1234"
$ws.Range("G48").Value = "Standar1234"
$ws.Range("H48").Value = "Standar1234"

$ws.Range("E49").Value = "MLWAG0215"
$ws.Range("F49").Value = "MLWAG0215"
$ws.Range("G49").Value = "MLWAG0215"
$ws.Range("H49").Value = "MLWAG0215"

$ws.Range("E50").Value = "TRAVELMADEZY"
$ws.Range("F50").Value = "TRAVELMADEZY"
$ws.Range("G50").Value = "TRAVELMADEZY"
$ws.Range("H50").Value = "TRAVELMADEZY"

# --- Column C: variable name -------------------------------------------
$ws.Range("C45").Value = "webUserName"
$ws.Range("C46").Value = "webPassword"
$ws.Range("C47").Value = "webUserName"
$ws.Range("C48").Value = "webPassword"
$ws.Range("C49").Value = "webUserName"
$ws.Range("C50").Value = "webPassword"

# --- Column A: project wildcard (matches row 44's "*") ----------------------
$ws.Range("A45").Value = "*"
$ws.Range("A46").Value = "*"
$ws.Range("A47").Value = "*"
$ws.Range("A48").Value = "*"
$ws.Range("A49").Value = "*"
$ws.Range("A50").Value = "*"

# --- Column D: generated variable-name formula, same pattern as D44 --------
$ws.Range("D45").Formula = '="@@."&A45&"."&B45&"."&C45&"@@"'
$ws.Range("D46").Formula = '="@@."&A46&"."&B46&"."&C46&"@@"'
$ws.Range("D47").Formula = '="@@."&A47&"."&B47&"."&C47&"@@"'
$ws.Range("D48").Formula = '="@@."&A48&"."&B48&"."&C48&"@@"'
$ws.Range("D49").Formula = '="@@."&A49&"."&B49&"."&C49&"@@"'
$ws.Range("D50").Formula = '="@@."&A50&"."&B50&"."&C50&"@@"'

# --- Style fix: column G on these rows should wrap text like E/F/H --------
$ws.Range("G45:G50").WrapText = $true

# --- Restore the column-A/B/C/D border + top-align styling that row 44
#     used, onto the newly populated cells (they already carried the
#     "empty row" style s=1, which matches; nothing further to do there).

# --- Update the view state: scroll position + active selection -------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C48").Select() | Out-Null
